# Apply column shift: drop the "Nominal 2023 GDP" column (old Q), shift the two
# "Sustainable Bonds Issued per GDP" columns (old R,S) left into Q,R, rename their
# headers, and remove the now-unused column S.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 43

# Capture old column R and S values (including headers) before overwriting,
# since we are shifting values left by one column.
$rVals = @{}
$sVals = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $rVals[$r] = $ws.Cells.Item($r, 18).Value2  # column R = 18
    $sVals[$r] = $ws.Cells.Item($r, 19).Value2  # column S = 19
}

# New headers for Q and R
$ws.Cells.Item(1, 17).Value2 = "Labelled_Sustainable_Bonds_per_GDP"
$ws.Cells.Item(1, 18).Value2 = "Unlabelled_Sustainable_Bonds_per_GDP"

# Shift the data rows: Q gets old R, R gets old S.
# Skip writing empty strings - the runtime deletes a cell outright when it is
# assigned an empty string, whereas an already-blank cell should just stay blank.
for ($r = 2; $r -le $lastRow; $r++) {
    if ("$($rVals[$r])" -ne "") {
        $ws.Cells.Item($r, 17).Value2 = $rVals[$r]
    }
    if ("$($sVals[$r])" -ne "") {
        $ws.Cells.Item($r, 18).Value2 = $sVals[$r]
    }
}

# Remove the now-obsolete column S entirely
$ws.Columns.Item(19).Delete()
